# Apply changes described by the diff to the Optical_Power worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows: OT (column E) values for rows 52 and 53 ---
$ws.Range('E52').Value = 'ICD31065552'
$ws.Range('E53').Value = 'ICD31075271'

# --- Append new rows 55-61 ---

# Force columns A (Caso) and B (F. De Reclamo) to be stored as plain text
# for the new rows, so values like "7289" and "9/24/2025" are not
# auto-converted by Excel into numbers/dates.
$ws.Range('A55:B61').NumberFormat = '@'

# Row 55
$ws.Range('A55').Value = '7289'
$ws.Range('B55').Value = '9/24/2025'
$ws.Range('C55').Value = 'ARANGUREN, JUAN F., DR. 1382'
$ws.Range('D55').Value = 6
$ws.Range('E55').Value = 'Pendiente ADM'
$ws.Range('F55').Value = 'Optical Power'
$ws.Range('G55').Value = 'Pendiente'
$ws.Range('H55').Value = 'Cable cortado'
$ws.Range('I55').Value = 1
$ws.Range('J55').Value = '{"direccionesNormalizadas": [{"altura": 1382, "cod_calle": 1094, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.451613", "y": "-34.618254"}, "direccion": "ARANGUREN, JUAN F., DR. 1382, CABA", "nombre_calle": "ARANGUREN, JUAN F., DR.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range('K55').Value = -58.451613
$ws.Range('L55').Value = -34.618254
$ws.Range('M55').Value = 'Boedo'
$ws.Range('N55').Value = 'Capital Sur'

# Row 56
$ws.Range('A56').Value = '7304'
$ws.Range('B56').Value = '9/25/2025'
$ws.Range('C56').Value = 'WASHINGTON 3926'
$ws.Range('D56').Value = 12
$ws.Range('E56').Value = 'ICD31075927'
$ws.Range('F56').Value = 'Optical Power'
$ws.Range('G56').Value = 'Pendiente'
$ws.Range('H56').Value = 'Cable en panza'
$ws.Range('I56').Value = 1
$ws.Range('J56').Value = '{"direccionesNormalizadas": [{"altura": 3926, "cod_calle": 24003, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.482836", "y": "-34.552950"}, "direccion": "WASHINGTON 3926, CABA", "nombre_calle": "WASHINGTON", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range('K56').Value = -58.482836
$ws.Range('L56').Value = -34.55295
$ws.Range('M56').Value = 'Saavedra'
$ws.Range('N56').Value = 'Capital Norte'

# Row 57
$ws.Range('A57').Value = '7312'
$ws.Range('B57').Value = '9/25/2025'
$ws.Range('C57').Value = 'BOGOTA 2902'
$ws.Range('D57').Value = 7
$ws.Range('E57').Value = 'ICD31076138'
$ws.Range('F57').Value = 'Optical Power'
$ws.Range('G57').Value = 'Pendiente'
$ws.Range('H57').Value = 'Tendido a baja altura y cortado'
$ws.Range('I57').Value = 1
$ws.Range('J57').Value = '{"direccionesNormalizadas": [{"altura": 2902, "cod_calle": 2090, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.471303", "y": "-34.627529"}, "direccion": "BOGOTA 2902, CABA", "nombre_calle": "BOGOTA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range('K57').Value = -58.471303
$ws.Range('L57').Value = -34.627529
$ws.Range('M57').Value = 'Devoto'
$ws.Range('N57').Value = 'Capital Norte'

# Row 58
$ws.Range('A58').Value = '7319'
$ws.Range('B58').Value = '9/26/2025'
$ws.Range('C58').Value = 'ALBERDI, JUAN BAUTISTA AV. 1723'
$ws.Range('D58').Value = 7
$ws.Range('E58').Value = 'Pendiente ADM'
$ws.Range('F58').Value = 'Optical Power'
$ws.Range('G58').Value = 'Pendiente'
$ws.Range('H58').Value = 'Tendido a baja altura'
$ws.Range('I58').Value = 1
$ws.Range('J58').Value = '{"direccionesNormalizadas": [{"altura": 1723, "cod_calle": 1033, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.453202", "y": "-34.628031"}, "direccion": "ALBERDI, JUAN BAUTISTA AV. 1723, CABA", "nombre_calle": "ALBERDI, JUAN BAUTISTA AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range('K58').Value = -58.453202
$ws.Range('L58').Value = -34.628031
$ws.Range('M58').Value = 'Boedo'
$ws.Range('N58').Value = 'Capital Sur'

# Row 59
$ws.Range('A59').Value = '7322'
$ws.Range('B59').Value = '9/26/2025'
$ws.Range('C59').Value = 'MORENO, JOSE MARIA AV. 864'
$ws.Range('D59').Value = 7
$ws.Range('E59').Value = 'Pendiente ADM'
$ws.Range('F59').Value = 'Optical Power'
$ws.Range('G59').Value = 'Pendiente'
$ws.Range('H59').Value = 'Cable en panza'
$ws.Range('I59').Value = 1
$ws.Range('J59').Value = '{"direccionesNormalizadas": [{"altura": 864, "cod_calle": 13128, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.433574", "y": "-34.628743"}, "direccion": "MORENO, JOSE MARIA AV. 864, CABA", "nombre_calle": "MORENO, JOSE MARIA AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range('K59').Value = -58.433574
$ws.Range('L59').Value = -34.628743
$ws.Range('M59').Value = 'Boedo'
$ws.Range('N59').Value = 'Capital Sur'

# Row 60
$ws.Range('A60').Value = '2219'
$ws.Range('B60').Value = '9/26/2025'
$ws.Range('C60').Value = 'CORDOBA AV. 3421'
$ws.Range('D60').Value = 2
$ws.Range('E60').Value = 'Pendiente ADM'
$ws.Range('F60').Value = 'Optical Power'
$ws.Range('G60').Value = 'Pendiente'
$ws.Range('H60').Value = 'Cable a baja altura '
$ws.Range('I60').Value = 1
$ws.Range('J60').Value = '{"direccionesNormalizadas": [{"altura": 3421, "cod_calle": 3165, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.414916", "y": "-34.597799"}, "direccion": "CORDOBA AV. 3421, CABA", "nombre_calle": "CORDOBA AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range('K60').Value = -58.414916
$ws.Range('L60').Value = -34.597799
$ws.Range('M60').Value = 'Almagro'
$ws.Range('N60').Value = 'Capital Sur'

# Row 61
$ws.Range('A61').Value = '2569'
$ws.Range('B61').Value = '9/26/2025'
$ws.Range('C61').Value = 'LAUTARO 1250'
$ws.Range('D61').Value = 7
$ws.Range('E61').Value = 'Pendiente ADM'
$ws.Range('F61').Value = 'Optical Power'
$ws.Range('G61').Value = 'Pendiente'
$ws.Range('H61').Value = 'CAbles a baja altura'
$ws.Range('I61').Value = 1
$ws.Range('J61').Value = '{"direccionesNormalizadas": [{"altura": 1250, "cod_calle": 12086, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.449719", "y": "-34.640108"}, "direccion": "LAUTARO 1250, CABA", "nombre_calle": "LAUTARO", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range('K61').Value = -58.449719
$ws.Range('L61').Value = -34.640108
$ws.Range('M61').Value = 'Boedo'
$ws.Range('N61').Value = 'Capital Sur'

# Reset style of the forced-text columns back to Normal/General so no
# residual explicit cell style is left behind (matches source formatting).
$ws.Range('A55:B61').Style = 'Normal'
